$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1960.875
$ws.Range("I15").Value = 1960.875
$ws.Range("K15").Value = 5882.625
$ws.Range("M15").Value = -5713.625

$ws.Range("H28").Value = 506.5
$ws.Range("I28").Value = 506.5
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 506.5
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -21.5
$ws.Range("N28").ClearContents()

$ws.Range("H40").Value = 5290.3335
$ws.Range("I40").Value = 3498.75
$ws.Range("K40").Value = 3498.75
$ws.Range("M40").Value = -3323.75

$ws.Range("H137").Value = 4358.75
$ws.Range("I137").Value = 4347
$ws.Range("J137").Value = 4362.6665
$ws.Range("K137").Value = 13041
$ws.Range("L137").Value = 13087.9995
$ws.Range("M137").Value = -10491
$ws.Range("N137").Value = -18187.9995

$ws.Range("H141").Value = 1498.5
$ws.Range("I141").Value = 1498.5
$ws.Range("K141").Value = 4495.5
$ws.Range("M141").Value = 684.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()

$ws.Range("H61").Value = 2900
$ws.Range("J61").Value = 5000
$ws.Range("L61").Value = 5000
$ws.Range("N61").Value = -5424

$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").ClearContents()

$ws.Range("H136").Value = 2900
$ws.Range("J136").Value = 5000
$ws.Range("L136").Value = 15000
$ws.Range("N136").Value = -20100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3859.2
$ws.Range("I105").Value = 3859.2
$ws.Range("K105").Value = 3859.2
$ws.Range("M105").Value = -2112.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3735.3333
$ws.Range("J62").Value = 3853
$ws.Range("L62").Value = 3853
$ws.Range("N62").Value = -5101

$ws.Range("H65").Value = 3735.3333
$ws.Range("J65").Value = 3853
$ws.Range("L65").Value = 19265
$ws.Range("N65").Value = -25505

$ws.Range("H68").Value = 78795
$ws.Range("J68").Value = 78795
$ws.Range("L68").Value = 78795
$ws.Range("N68").Value = -80293

$ws.Range("H71").Value = 78795
$ws.Range("J71").Value = 78795
$ws.Range("L71").Value = 236385
$ws.Range("N71").Value = -243873

$ws.Range("H74").Value = 72814
$ws.Range("J74").Value = 72814
$ws.Range("L74").Value = 72814
$ws.Range("N74").Value = -74562

$ws.Range("H77").Value = 72814
$ws.Range("J77").Value = 72814
$ws.Range("L77").Value = 218442
$ws.Range("N77").Value = -227178

$ws.Range("H86").Value = 4145.6313
$ws.Range("I86").Value = 2052.889
$ws.Range("K86").Value = 2052.889
$ws.Range("M86").Value = -929.8890000000001

$ws.Range("H89").Value = 4145.6313
$ws.Range("I89").Value = 2052.889
$ws.Range("K89").Value = 10264.445
$ws.Range("M89").Value = -4648.445

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 3000
$ws.Range("I14").Value = 3000
$ws.Range("K14").Value = 9000
$ws.Range("M14").Value = -8827

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 3
$ws.Range("I11").Value = 3
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 3
$ws.Range("L11").ClearContents()
$ws.Range("M11").Value = 136
$ws.Range("N11").ClearContents()

$ws.Range("H80").Value = 9000
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 9000
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 9000
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -10996

$ws.Range("H83").Value = 9000
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 9000
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 45000
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -54984

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1486.875
$ws.Range("I22").Value = 849.25
$ws.Range("J22").Value = 2124.5
$ws.Range("K22").Value = 849.25
$ws.Range("L22").Value = 2124.5
$ws.Range("M22").Value = -554.25
$ws.Range("N22").Value = -2714.5

$ws.Range("H27").Value = 1486.875
$ws.Range("I27").Value = 849.25
$ws.Range("J27").Value = 2124.5
$ws.Range("K27").Value = 849.25
$ws.Range("L27").Value = 2124.5
$ws.Range("M27").Value = -742.25
$ws.Range("N27").Value = -2338.5

$ws.Range("H68").Value = 1850
$ws.Range("I68").Value = 1850
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 1850
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -1101
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 1850
$ws.Range("I71").Value = 1850
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 9250
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -5506
$ws.Range("N71").ClearContents()

$ws.Range("H93").Value = 2747.75
$ws.Range("I93").Value = 3496.5
$ws.Range("K93").Value = 3496.5
$ws.Range("M93").Value = -2248.5

$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").ClearContents()

$ws.Range("H17").Value = 3800
$ws.Range("I17").Value = 3800
$ws.Range("K17").Value = 3800
$ws.Range("M17").Value = -3628

$ws.Range("H81").Value = 16813.182
$ws.Range("I81").Value = 17494.5
$ws.Range("J81").Value = 10000
$ws.Range("K81").Value = 34989
$ws.Range("L81").Value = 20000
$ws.Range("M81").Value = -33928
$ws.Range("N81").Value = -22122

$ws.Range("H84").Value = 16813.182
$ws.Range("I84").Value = 17494.5
$ws.Range("J84").Value = 10000
$ws.Range("K84").Value = 174945
$ws.Range("L84").Value = 100000
$ws.Range("M84").Value = -169641
$ws.Range("N84").Value = -110608

$ws.Range("H104").Value = 15000
$ws.Range("J104").Value = 15000
$ws.Range("L104").Value = 15000
$ws.Range("N104").Value = -21988
